$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.676.62"
$ws.Range("E2").Value = "  +1.43%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.160.06"
$ws.Range("E3").Value = "  +0.88%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "529.87"
$ws.Range("E5").Value = "  -0.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.93"
$ws.Range("E6").Value = "  +0.90%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("E8").Value = "  +16.75%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.31"
$ws.Range("E9").Value = "  +0.24%  "

$ws.Range("E10").Value = "  +5.47%  "

$ws.Range("E11").Value = "  +4.13%  "

$ws.Range("E12").Value = "  +3.48%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.704.86"
$ws.Range("E13").Value = "  +0.95%  "

$ws.Range("E14").Value = "  +1.26%  "

$ws.Range("E15").Value = "  +4.84%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.728.07"
$ws.Range("E16").Value = "  +1.27%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.26"
$ws.Range("E17").Value = "  +3.72%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.132.17"
$ws.Range("E18").Value = "  +0.11%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.00"
$ws.Range("E19").Value = "  +2.23%  "

$ws.Range("E20").Value = "  +0.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "376.46"
$ws.Range("E21").Value = "  +4.64%  "

$ws.Range("E22").Value = "  +1.89%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.12%  "

$ws.Range("E24").Value = "  +5.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.74"
$ws.Range("E25").Value = "  +1.16%  "

$ws.Range("E26").Value = "  +0.48%  "

$ws.Range("E27").Value = "  -0.13%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.31"
$ws.Range("E28").Value = "  +13.83%  "

$ws.Range("E29").Value = "  -0.62%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.34"
$ws.Range("E30").Value = "  +4.34%  "

$ws.Range("E31").Value = "  +0.45%  "

$ws.Range("E32").Value = "  -0.61%  "

$ws.Range("E33").Value = "  +0.25%  "

$ws.Range("E34").Value = "  +1.12%  "

$ws.Range("E35").Value = "  +3.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "158.10"
$ws.Range("E36").Value = "  -0.23%  "

$ws.Range("E37").Value = "  +4.97%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.00"
$ws.Range("E38").Value = "  -3.39%  "

$ws.Range("E39").Value = "  +1.86%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0690"
$ws.Range("E40").Value = "  +2.64%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.647.17"
$ws.Range("E41").Value = "  +5.85%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.28"
$ws.Range("E42").Value = "  +7.16%  "

$ws.Range("E43").Value = "  +3.34%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.15"
$ws.Range("E44").Value = "  +3.78%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0288"
$ws.Range("E45").Value = "  +7.82%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  -0.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.201.91"
$ws.Range("E47").Value = "  +0.84%  "

$ws.Range("E48").Value = "  +14.76%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.21"
$ws.Range("E49").Value = "  +2.17%  "

$ws.Range("E50").Value = "  -1.22%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.06"
$ws.Range("E51").Value = "  +1.43%  "
